# Added RMA Test Cases: SO To inspection order SO to RMA Receipt
#
# The "RMA Details Maintenance Grid" sheet holds one RMA test-data group
# (RMA number / shipper-line number / Salesforce record id) per data row.
# This commit swaps the most recent group (RMA-C084-*) out for a freshly
# generated group (RMA-CL8I-*) across the three data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2
$ws.Range("E2").Value = "RMA-CL8I-001"
$ws.Range("F2").Value = "RMA-CL8I-1-1"
$ws.Range("J2").Value = "a7s5f000000xK54AAE"

# Row 3
$ws.Range("E3").Value = "RMA-CL8I-002"
$ws.Range("F3").Value = "RMA-CL8I-1-2"
$ws.Range("J3").Value = "a7s5f000000xK55AAE"

# Row 4
$ws.Range("E4").Value = "RMA-CL8I-003"
$ws.Range("F4").Value = "RMA-CL8I-1-3"
$ws.Range("J4").Value = "a7s5f000000xK56AAE"
